$d = $word.ActiveDocument

# 1) Split the "Critério" run: insert a manual line break after the colon
#    sentence, before "Média final = ..."
$d.Content.Find.Execute(
    "A média final será calculada segundo a equação abaixo: Média final",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A média final será calculada segundo a equação abaixo: ^lMédia final",
    2)

# 2) Split the bibliography run into four lines separated by manual line
#    breaks.
$d.Content.Find.Execute(
    "1. Nelson, D.L., Cox, M.M. Princípios de bioquímica de Lehninger. ArtmedEditora, 2022.2. Segel, I.H. Bioquímica Teoria e Problemas, São Paulo: Livros técnicos e Científicos Editora S.A, 1979. 3. Artigos e revisões da literatura ou outra bibliografia indicada no cronograma anual da disciplina.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "1. Nelson, D.L., Cox, M.M. Princípios de bioquímica de Lehninger. Artmed^lEditora, 2022.^l2. Segel, I.H. Bioquímica Teoria e Problemas, São Paulo: Livros técnicos e Científicos Editora S.A, 1979.^l 3. Artigos e revisões da literatura ou outra bibliografia indicada no cronograma anual da disciplina.",
    2)
